$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.958.22'
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").Value = '1.641.59'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("E4").Value = '  +0.08%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '212.64'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.59%  '
$ws.Range("E6").Value = '  +0.42%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.05%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '23.48'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +1.77%  '
$ws.Range("E9").Value = '  -1.67%  '
$ws.Range("E11").Value = '  +2.55%  '
$ws.Range("D12").Value = '1.873.61'
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("D13").Value = '1.640.30'
$ws.Range("E13").Value = '  +0.15%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '4.09'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.05%  '
$ws.Range("E15").Value = '  +2.01%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '65.51'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("D17").Value = '27.948.69'
$ws.Range("E17").Value = '  +1.28%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '233.09'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.65%  '
$ws.Range("E19").Value = '  +0.56%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '7.60'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.71%  '
$ws.Range("E21").Value = '  +0.11%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '10.58'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.48%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '4.38'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.23%  '
$ws.Range("E24").Value = '  -1.52%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '153.15'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +2.70%  '
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("E27").Value = '  +0.34%  '
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  +0.46%  '
$ws.Range("E32").Value = '  +4.02%  '
$ws.Range("E33").Value = '  +0.51%  '
$ws.Range("D34").Value = '1.407.06'
$ws.Range("E34").Value = '  -3.81%  '
$ws.Range("E35").Value = '  +2.40%  '
$ws.Range("E36").Value = '  +1.77%  '
$ws.Range("E37").Value = '  +1.61%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.564'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +1.20%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.880'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +0.18%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.927'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +0.31%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '1.03'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +1.18%  '
$ws.Range("E42").Value = '  +0.05%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '67.24'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.43%  '
$ws.Range("E44").Value = '  +6.24%  '
$ws.Range("E45").Value = '  +2.48%  '
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("D47").Value = '1.782.90'
$ws.Range("E47").Value = '  +0.65%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '87.85'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.56%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.100'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("E50").Value = '  +0.36%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '7.61'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.49%  '
